$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
foreach ($shp in $s.Shapes) {
    if ($shp.HasTable) {
        $tbl = $shp.Table
        $tbl.ApplyStyle("{10890D0B-0F45-4CAE-8E81-8D9B11325874}")
    }
}
